# AutoCommit_10 июля 2024 г. 15:28:43_SibNout2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in assessment marks (value 5) for the relevant cells
$ws.Range("H6").Value = 5

$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 5
$ws.Range("E9").Value = 5

$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 5
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5

$ws.Range("D14").Value = 5
$ws.Range("H14").Value = 5

$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 5
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 5

$ws.Range("E23").Value = 5
$ws.Range("H23").Value = 5

$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 5
$ws.Range("H24").Value = 5

$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 5
$ws.Range("I26").Value = 5

$ws.Range("H27").Value = 5

$ws.Range("G29").Value = 5
$ws.Range("H29").Value = 5
$ws.Range("I29").Value = 5

# Update the active selection in the frozen-pane view to C9
$ws.Range("C9").Select()
